# Poster_RS 07.29 — "Updated to Reflect Abstract Submission"
# Re-positions / re-sizes several boxes & images, switches body-text
# paragraphs from left/center aligned to justified, and re-crops one
# of the method photos.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# ---- helpers -----------------------------------------------------
function Get-ShapeById($slide, $id) {
    for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
        $sh = $slide.Shapes.Item($i)
        if ($sh.Id -eq $id) { return $sh }
    }
    return $null
}

function EMU-ToPoints($emu) {
    return $emu / 914400.0 * 72.0
}

function Set-ShapePos($shape, $x, $y) {
    $shape.Left = EMU-ToPoints $x
    $shape.Top  = EMU-ToPoints $y
}

function Set-ShapeSize($shape, $cx, $cy) {
    $shape.Width  = EMU-ToPoints $cx
    $shape.Height = EMU-ToPoints $cy
}

function Set-Justify($shape) {
    $shape.TextFrame.TextRange.ParagraphFormat.Alignment = 4   # ppAlignJustify
}

# ---- 1. "Hypothesis" heading pill (Rectangle: Rounded Corners 84) --
$sh = Get-ShapeById $s 85
Set-ShapePos  $sh 689679 10284139
Set-ShapeSize $sh 15333105 806098

# ---- 2. "Design" heading pill (Rectangle: Rounded Corners 28) ------
$sh = Get-ShapeById $s 29
Set-ShapePos  $sh 16936458 4319973
Set-ShapeSize $sh 31119871 813098

# ---- 3. Intro paragraph (TextBox 59) — keep marL/marR, add justify -
$sh = Get-ShapeById $s 60
$sh.TextFrame.TextRange.Paragraphs(1, 1).ParagraphFormat.Alignment = 4

# ---- 4. "S-ʃ sounds..." box (TextBox 72) ---------------------------
$sh = Get-ShapeById $s 73
Set-ShapePos $sh 17299554 5857363
Set-Justify  $sh

# ---- 5. "In this study..." box (TextBox 76) — move only -----------
$sh = Get-ShapeById $s 77
Set-ShapePos $sh 20608776 5179319

# ---- 6. References list (TextBox 81) — justify the Luthra entry ---
$sh = Get-ShapeById $s 82
$sh.TextFrame.TextRange.Paragraphs(4, 1).ParagraphFormat.Alignment = 4

# ---- 7. "Our critical stimuli..." (TextBox 82) ---------------------
$sh = Get-ShapeById $s 83
Set-Justify $sh

# ---- 8. "All recordings were processed..." (TextBox 83) -----------
$sh = Get-ShapeById $s 84
Set-Justify $sh

# ---- 9. "Participants will be instructed..." (TextBox 68) ---------
$sh = Get-ShapeById $s 69
Set-Justify $sh

# ---- 10. Method photo (Picture 98) — move only ---------------------
$sh = Get-ShapeById $s 99
Set-ShapePos $sh 41627051 11632440

# ---- 11. "The word pairings..." (TextBox 26) -----------------------
$sh = Get-ShapeById $s 27
Set-Justify $sh

# ---- 12. "These words were split..." (TextBox 61) ------------------
$sh = Get-ShapeById $s 62
Set-Justify $sh

# ---- 13. "Talker A and Talker B..." (TextBox 45) -------------------
$sh = Get-ShapeById $s 46
Set-Justify $sh

# ---- 14. "Each experiment consists..." (TextBox 58) ----------------
$sh = Get-ShapeById $s 59
Set-Justify $sh

# ---- 15. "A listener's perceptual boundary..." (TextBox 62) -------
$sh = Get-ShapeById $s 63
Set-Justify $sh

# ---- 16. "After the Exposure Phase..." (TextBox 77) ----------------
$sh = Get-ShapeById $s 78
Set-Justify $sh

# ---- 17. Acknowledgements funding blurb (TextBox 79) --------------
$sh = Get-ShapeById $s 80
Set-ShapeSize $sh 8583002 1477328
Set-Justify   $sh

# ---- 18. Acknowledgements thank-you blurb (TextBox 80) ------------
$sh = Get-ShapeById $s 81
Set-Justify $sh

# ---- 19. "Predictions" heading pill (Rectangle: Rounded Corners 39) -
$sh = Get-ShapeById $s 40
Set-ShapePos  $sh 834487 14781891
Set-ShapeSize $sh 15228809 806099

# ---- 20. Predictions caption (TextBox 88) — move only --------------
$sh = Get-ShapeById $s 89
Set-ShapePos $sh 795700 26385556

# ---- 21. Predictions chart photo (Picture 14) — re-crop + reposition
$sh = Get-ShapeById $s 15
$sh.PictureFormat.CropBottom = 37.175
Set-ShapePos  $sh 2834036 15587990
Set-ShapeSize $sh 11300888 5104893

# ---- 22. Second predictions chart photo (Picture 9) — reposition ---
$sh = Get-ShapeById $s 10
Set-ShapePos  $sh 2882423 20656381
Set-ShapeSize $sh 11300888 5749229
